$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.332116365432739
$ws.Range("B1").Value = 1.52305018901825
$ws.Range("C1").Value = 6.872737407684326
$ws.Range("D1").Value = 1.939961433410645
$ws.Range("E1").Value = 0.879861056804657
